# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps written by the handback report
# generator.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date (shared value also shown on de-de!H2)
$overview.Range("G2").Value = "2016-09-01 19:15:58"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-09-01 19:15:53"
$zhcn.Range("K2").Value = "2016-09-01 19:16:23"

# de-de: Correspond Handoff Datetime (mirrors Overview!G2) / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-09-01 19:15:58"
$dede.Range("K2").Value = "2016-09-01 19:16:30"
